# Saldo.xlsx edit
# - Remove the account "001761119" (BLUEMETRIX, 454.16) row.
# - Add a new account "005993550" (ALESSANDRA, 288.1) row, keeping the
#   sheet's descending-balance sort order (it belongs right after the
#   "005547703" / SILVIA / 290.74 row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 142 holds 001761119 / BLUEMETRIX / 454.16 - delete it entirely,
# shifting all following rows up by one.
$ws.Rows.Item(142).Delete()

# After that deletion, the SILVIA / 290.74 row is now row 181, and the
# row that follows it (ALESSANDRA / 284.74) is now row 182. Insert a
# fresh blank row at 182 so the new entry lands between them.
$ws.Rows.Item(182).Insert()

# Account numbers are stored as text (leading zeros matter), so force
# the cell to Text before writing the value, then clear the formatting
# again so no stray number-format style is left behind on the cell.
$ws.Cells.Item(182, 1).NumberFormat = "@"
$ws.Cells.Item(182, 1).Value = "005993550"
$ws.Cells.Item(182, 1).ClearFormats()

$ws.Cells.Item(182, 2).Value = "ALESSANDRA"
$ws.Cells.Item(182, 3).Value = 288.1
